# Replace the working set of sequences: new images, words, and categories
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "house/house016.jpg"
$ws.Range("C2").Value = "haken"
$ws.Range("D2").Value = "house"
$ws.Range("B3").Value = "house/house025.jpg"
$ws.Range("C3").Value = "mieten"
$ws.Range("D3").Value = "house"
$ws.Range("B4").Value = "house/house017.jpg"
$ws.Range("C4").Value = "gründen"
$ws.Range("D4").Value = "house"
$ws.Range("B5").Value = "dog/dog020.jpg"
$ws.Range("C5").Value = "pflegen"
$ws.Range("D5").Value = "dog"
$ws.Range("B6").Value = "dog/dog031.jpg"
$ws.Range("C6").Value = "lehnen"
$ws.Range("D6").Value = "dog"
$ws.Range("B7").Value = "dog/dog017.jpg"
$ws.Range("C7").Value = "opfern"
$ws.Range("D7").Value = "dog"
$ws.Range("B8").Value = "dog/dog001.jpg"
$ws.Range("C8").Value = "jubeln"
$ws.Range("D8").Value = "dog"
$ws.Range("B9").Value = "dog/dog004.jpg"
$ws.Range("C9").Value = "antun"
$ws.Range("D9").Value = "dog"
$ws.Range("B10").Value = "house/house013.jpg"
$ws.Range("C10").Value = "hoffen"
$ws.Range("D10").Value = "house"
$ws.Range("B11").Value = "house/house000.jpg"
$ws.Range("C11").Value = "dauern"
$ws.Range("D11").Value = "house"
$ws.Range("B12").Value = "house/house030.jpg"
$ws.Range("C12").Value = "rasen"
$ws.Range("D12").Value = "house"
$ws.Range("B13").Value = "dog/dog015.jpg"
$ws.Range("C13").Value = "wiegen"
$ws.Range("D13").Value = "dog"
$ws.Range("B14").Value = "house/house006.jpg"
$ws.Range("C14").Value = "kaufen"
$ws.Range("D14").Value = "house"
$ws.Range("B15").Value = "dog/dog005.jpg"
$ws.Range("C15").Value = "stärken"
$ws.Range("D15").Value = "dog"
$ws.Range("B16").Value = "house/house009.jpg"
$ws.Range("C16").Value = "stechen"
$ws.Range("D16").Value = "house"
$ws.Range("B17").Value = "house/house028.jpg"
$ws.Range("C17").Value = "scheitern"
$ws.Range("D17").Value = "house"
$ws.Range("B18").Value = "house/house022.jpg"
$ws.Range("C18").Value = "drehen"
$ws.Range("D18").Value = "house"
$ws.Range("B19").Value = "house/house007.jpg"
$ws.Range("C19").Value = "töten"
$ws.Range("D19").Value = "house"
$ws.Range("B20").Value = "dog/dog009.jpg"
$ws.Range("C20").Value = "formen"
$ws.Range("D20").Value = "dog"
$ws.Range("B21").Value = "house/house019.jpg"
$ws.Range("C21").Value = "regnen"
$ws.Range("D21").Value = "house"
$ws.Range("B22").Value = "dog/dog026.jpg"
$ws.Range("C22").Value = "schätzen"
$ws.Range("D22").Value = "dog"
$ws.Range("B23").Value = "house/house021.jpg"
$ws.Range("C23").Value = "wenden"
$ws.Range("D23").Value = "house"
$ws.Range("B24").Value = "house/house005.jpg"
$ws.Range("C24").Value = "nehmen"
$ws.Range("D24").Value = "house"
$ws.Range("B25").Value = "dog/dog007.jpg"
$ws.Range("C25").Value = "fliegen"
$ws.Range("D25").Value = "dog"
$ws.Range("B26").Value = "dog/dog012.jpg"
$ws.Range("C26").Value = "backen"
$ws.Range("D26").Value = "dog"
$ws.Range("B27").Value = "house/house027.jpg"
$ws.Range("C27").Value = "enden"
$ws.Range("D27").Value = "house"
$ws.Range("B28").Value = "house/house026.jpg"
$ws.Range("C28").Value = "rücken"
$ws.Range("D28").Value = "house"
$ws.Range("B29").Value = "dog/dog006.jpg"
$ws.Range("C29").Value = "strahlen"
$ws.Range("D29").Value = "dog"
$ws.Range("B30").Value = "dog/dog025.jpg"
$ws.Range("C30").Value = "posten"
$ws.Range("D30").Value = "dog"
$ws.Range("B31").Value = "dog/dog028.jpg"
$ws.Range("C31").Value = "segeln"
$ws.Range("D31").Value = "dog"
$ws.Range("B32").Value = "dog/dog023.jpg"
$ws.Range("C32").Value = "husten"
$ws.Range("D32").Value = "dog"
$ws.Range("B33").Value = "dog/dog008.jpg"
$ws.Range("C33").Value = "krachen"
$ws.Range("D33").Value = "dog"
